# Add 2022 mortality data
# - Reorder a few comorbidity rows ("Preterm birth" now comes before
#   "Cardiac insufficiency", "Diabetes" now comes before "Asthma")
# - Fix a few typos in comorbidity names
# - Update several counts in column B to reflect 2022 mortality data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (comorbidity names) and column B (counts) for rows 2-20.
$ws.Range("A2").Value = "None reported"
$ws.Range("B2").Value = 140

$ws.Range("A3").Value = "Multiple comorbidities"
$ws.Range("B3").Value = 68

$ws.Range("A4").Value = "Neurological disease"
$ws.Range("B4").Value = 49

$ws.Range("A5").Value = "Oncological disease"
$ws.Range("B5").Value = 33

$ws.Range("A6").Value = "No comorbidities"
$ws.Range("B6").Value = 25

$ws.Range("A7").Value = "Preterm birth"
$ws.Range("B7").Value = 14

$ws.Range("A8").Value = "Cardiac insufficiency"
$ws.Range("B8").Value = 12

$ws.Range("A9").Value = "Congenital or acquired immunosuppression"
$ws.Range("B9").Value = 11

$ws.Range("A10").Value = "Diabetes"
$ws.Range("B10").Value = 6

$ws.Range("A11").Value = "Asthma"
$ws.Range("B11").Value = 5

$ws.Range("A12").Value = "Obesity"
$ws.Range("B12").Value = 5

$ws.Range("A13").Value = "Chronic renal insufficiency"
$ws.Range("B13").Value = 4

$ws.Range("A14").Value = "Low body weight"
$ws.Range("B14").Value = 3

$ws.Range("A15").Value = "Dialisis"
$ws.Range("B15").Value = 2

$ws.Range("A16").Value = "Previous community-acquired pneumonia"
$ws.Range("B16").Value = 2

$ws.Range("A17").Value = "Tuberculosis"
$ws.Range("B17").Value = 2

$ws.Range("A18").Value = "Arterial hypertension"
$ws.Range("B18").Value = 1

$ws.Range("A19").Value = "Chronic Obstructive Pulmonary Disease"
$ws.Range("B19").Value = 1

$ws.Range("A20").Value = "Previous bronchiolitis"
$ws.Range("B20").Value = 1
